$wb = $excel.ActiveWorkbook

# --- Group "18" submission: add three new names to sheet "15" (5th tab) ---
$ws = $wb.Worksheets.Item(5)
$ws.Activate() | Out-Null

# Enter the new names in the same order the source workbook's shared-string
# table lists them (Lior Malik, Astar Avraham, Noa Malka) even though they
# land in rows 8, 9, 7 respectively.
$ws.Range("A8").Value = "Lior Malik"
$ws.Range("A9").Value = "Astar Avraham"
$ws.Range("A7").Value = "Noa Malka"

# Leave the selection on the first of the newly entered cells, matching the
# workbook's last saved cursor position.
$ws.Range("A7").Select() | Out-Null

# --- Default workbook font: Calibri -> Arial ---
$wb.Styles.Item("Normal").Font.Name = "Arial"
